# Fruta / hortaliza, semanal
# Insert two new weekly price rows for "Alcachofa" (Terminal Hortofrutícola
# Agro Chillán) right after the existing row 44, pushing the previously
# existing rows 45-62 down to 47-64.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at position 45 (existing rows 45..62 shift to 47..64)
$ws.Rows("45:46").Insert()

# --- New row 45 ---
$ws.Range("A45").Value = 7
$ws.Range("B45").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C45").Value = "Ñuble"
$ws.Range("D45").Value = 44813
$ws.Range("D45").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E45").Value = 16
$ws.Range("F45").Value = 100112013
$ws.Range("G45").Value = "Alcachofa"
$ws.Range("H45").Value = "Argentina(o)"
$ws.Range("I45").Value = "Primera"
$ws.Range("J45").Value = 80
$ws.Range("K45").Value = 13000
$ws.Range("L45").Value = 14000
$ws.Range("M45").Value = 13500
$ws.Range("N45").Value = "`$/caja 50 unidades"
$ws.Range("O45").Value = "Provincia de Limarí"
$ws.Range("P45").Value = 270
$ws.Range("Q45").Value = 50
$ws.Range("R45").Value = "Hortaliza"

# --- New row 46 ---
$ws.Range("A46").Value = 7
$ws.Range("B46").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C46").Value = "Ñuble"
$ws.Range("D46").Value = 44813
$ws.Range("D46").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E46").Value = 16
$ws.Range("F46").Value = 100112013
$ws.Range("G46").Value = "Alcachofa"
$ws.Range("H46").Value = "Madrigal"
$ws.Range("I46").Value = "Primera"
$ws.Range("J46").Value = 100
$ws.Range("K46").Value = 11000
$ws.Range("L46").Value = 12000
$ws.Range("M46").Value = 11500
$ws.Range("N46").Value = "`$/caja 40 unidades"
$ws.Range("O46").Value = "Provincia de Limarí"
$ws.Range("P46").Value = 288
$ws.Range("Q46").Value = 40
$ws.Range("R46").Value = "Hortaliza"
